$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The customer row for phone "09876543" (row 46) is being split in two:
#  - row 46 keeps its birthday/points but the phone becomes the numeric 9876543
#    (the leading zero is dropped once it's stored as a number)
#  - a new row 47 is inserted for the original text phone "09876543" with
#    points reset to 0 and no birthday on file yet.

# Insert a fresh blank row directly below the existing row 46.
$ws.Rows.Item(47).Insert()

# New row 47: phone kept as text (leading zero preserved via quote-prefix),
# birthday blank, points 0.
$ws.Range("A47").Value = "'09876543"
$ws.Range("A47").Style = "Normal"

$ws.Range("B47").Value = "'"
$ws.Range("B47").Style = "Normal"

$ws.Range("C47").Value = 0

# Existing row 46: phone becomes a plain number (0.00 points, unchanged).
$ws.Range("A46").Value = 9876543
